# Update "latest output (run 24)" optimisation_result.xlsx
#
# Sheet "Schedule": rows 5 & 6 are removed (only 3 data rows remain),
#   and the values of rows 2-4 are refreshed with new optimisation output.
# Sheet "Detailed": many Price (col B) values were refreshed with new
#   forecast/historical data, some Type (col C) values flip between
#   "forecast" and "historical", and several Pump_Status (col E) flags
#   flip between "ON" and "OFF".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Schedule"
# ---------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Remove the two trailing rows (6 then 5) first, from bottom to top, so
# row numbers of the rows we keep (2-4) do not shift while deleting.
$schedule.Rows.Item(6).Delete()
$schedule.Rows.Item(5).Delete()

# Refresh the remaining three data rows with the new run's values.
$schedule.Range("A2").Value = 46037
$schedule.Range("B2").Value = 46037.66666666666
$schedule.Range("C2").Value = 16
$schedule.Range("D2").Value = 60.48
$schedule.Range("E2").Value = 1670.962254
$schedule.Range("F2").Value = 27.62834414682539

$schedule.Range("A3").Value = 46038.33333333334
$schedule.Range("B3").Value = 46038.66666666666
$schedule.Range("C3").Value = 8
$schedule.Range("D3").Value = 30.24
$schedule.Range("E3").Value = 620.6932777500001
$schedule.Range("F3").Value = 20.52557135416667

$schedule.Range("A4").Value = 46038.83333333334
$schedule.Range("B4").Value = 46039
# C4 (4) and D4 (15.12) are unchanged.
$schedule.Range("E4").Value = 331.7085135
$schedule.Range("F4").Value = 21.93839375

# ---------------------------------------------------------------------
# Sheet 2: "Detailed"
# ---------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

# Pump_Status flips OFF -> ON
$detailed.Range("E2").Value = "ON"
$detailed.Range("E3").Value = "ON"
$detailed.Range("E12").Value = "ON"
$detailed.Range("E13").Value = "ON"
$detailed.Range("E14").Value = "ON"
$detailed.Range("E15").Value = "ON"
$detailed.Range("E66").Value = "ON"
$detailed.Range("E67").Value = "ON"

# Pump_Status flips ON -> OFF
$detailed.Range("E44").Value = "OFF"
$detailed.Range("E45").Value = "OFF"
$detailed.Range("E46").Value = "OFF"
$detailed.Range("E47").Value = "OFF"
$detailed.Range("E48").Value = "OFF"
$detailed.Range("E49").Value = "OFF"
$detailed.Range("E50").Value = "OFF"
$detailed.Range("E51").Value = "OFF"

# Price (col B) refreshed values
$detailed.Range("B45").Value = 71.40000000000001
$detailed.Range("B46").Value = 65.0001
$detailed.Range("B47").Value = 65
$detailed.Range("B48").Value = 85.95
$detailed.Range("B49").Value = 84.7901
$detailed.Range("B50").Value = 81.05549000000001
$detailed.Range("B52").Value = 78
$detailed.Range("B53").Value = 78
$detailed.Range("B54").Value = 78
$detailed.Range("B55").Value = 69.42238
$detailed.Range("B56").Value = 60.46039
$detailed.Range("B57").Value = 59.40249
$detailed.Range("B58").Value = 59.39334
$detailed.Range("B59").Value = 79.95005
$detailed.Range("B60").Value = 79.95
$detailed.Range("B64").Value = 50.62891
$detailed.Range("B65").Value = 56.98
$detailed.Range("B67").Value = 56.97989
$detailed.Range("B68").Value = 47.31837
$detailed.Range("B70").Value = 45.92104
$detailed.Range("B71").Value = 45.97441
$detailed.Range("B72").Value = 36.06028
$detailed.Range("B75").Value = 40.54
$detailed.Range("B79").Value = 36.0601
$detailed.Range("B80").Value = 28.73596
$detailed.Range("B81").Value = 25.59822
$detailed.Range("B82").Value = 18.17021
$detailed.Range("B83").Value = 8.64973
$detailed.Range("B84").Value = -7.981
$detailed.Range("B85").Value = -6.80121
$detailed.Range("B86").Value = -6
$detailed.Range("B87").Value = -3.03118
$detailed.Range("B88").Value = -3.04997
$detailed.Range("B89").Value = 22.01959
$detailed.Range("B90").Value = 22.01959
$detailed.Range("B91").Value = 29.85322
$detailed.Range("B92").Value = 0.85459
$detailed.Range("B93").Value = 57.78152
$detailed.Range("B94").Value = 57.04922
$detailed.Range("B95").Value = 57.45801
$detailed.Range("B96").Value = 58.21771
$detailed.Range("B97").Value = 56.98

# Type (col C) flips forecast -> historical
$detailed.Range("C47").Value = "historical"
$detailed.Range("C48").Value = "historical"
$detailed.Range("C49").Value = "historical"
